$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-39, 45-51: refreshed Price (D) and Volume(1h) (E) figures.
# Numeric-looking Price strings are written via a leading apostrophe (forces
# Excel to store them as text instead of auto-converting to a Number), then the
# cell style is reset to "Normal" so no stray quote-prefix / text-format style
# lingers (matches the original plain inline-string cells).

$ws.Range("D2").Value = "62.377.10"
$ws.Range("E2").Value = "  -2.05%  "

$ws.Range("D3").Value = "2.434.57"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Formula = "'569.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "

$ws.Range("D6").Formula = "'143.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.65%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").Formula = "'0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.38%  "

$ws.Range("D9").Value = "2.430.62"
$ws.Range("E9").Value = "  -2.12%  "

$ws.Range("D10").Formula = "'0.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.72%  "

$ws.Range("D12").Formula = "'5.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.09%  "

$ws.Range("D13").Formula = "'0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.20%  "

$ws.Range("D14").Formula = "'26.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.35%  "

$ws.Range("D15").Formula = "'0.0000174"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.01%  "

$ws.Range("D16").Value = "2.878.50"
$ws.Range("E16").Value = "  -0.92%  "

$ws.Range("D17").Value = "62.248.99"
$ws.Range("E17").Value = "  -1.94%  "

$ws.Range("D18").Value = "2.434.25"
$ws.Range("E18").Value = "  -1.98%  "

$ws.Range("D19").Formula = "'11.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.36%  "

$ws.Range("D20").Formula = "'7.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.19%  "

$ws.Range("D21").Formula = "'324.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.39%  "

$ws.Range("D22").Formula = "'4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.74%  "

$ws.Range("D23").Formula = "'2.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.33%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").Formula = "'65.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.53%  "

$ws.Range("D26").Formula = "'619.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.44%  "

$ws.Range("D27").Formula = "'9.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.11%  "

$ws.Range("D28").Value = "0.0₃0961"
$ws.Range("E28").Value = "  -9.74%  "

$ws.Range("D29").Value = "2.553.73"
$ws.Range("E29").Value = "  -1.66%  "

$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").Formula = "'1.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.51%  "

$ws.Range("D32").Formula = "'8.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.62%  "

$ws.Range("D33").Formula = "'1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.32%  "

$ws.Range("D34").Formula = "'0.135"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.92%  "

$ws.Range("D35").Formula = "'5.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.94%  "

$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").Formula = "'1.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.69%  "

$ws.Range("D38").Formula = "'0.375"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.03%  "

$ws.Range("D39").Formula = "'18.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.60%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Formula = "'147.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Formula = "'5.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.80%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Formula = "'1.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.25%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Formula = "'42.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.64%  "

$ws.Range("D45").Formula = "'2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.49%  "

$ws.Range("D46").Formula = "'144.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.92%  "

$ws.Range("D47").Formula = "'3.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.26%  "

$ws.Range("D48").Formula = "'0.0523"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.36%  "

$ws.Range("D49").Formula = "'20.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.04%  "

$ws.Range("D50").Formula = "'0.594"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.85%  "

$ws.Range("D51").Formula = "'0.0229"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.84%  "

